$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. "...using Visual Studio 2019, including the following steps" ->
#    "...using Visual Studio 2019. The process includes the following steps"
Replace-Text " using Visual Studio 2019, including the following steps" " using Visual Studio 2019. The process includes the following steps"

# 2. "Imported the necessary modules" -> "Imported the necessary modules/libraries"
Replace-Text "Imported the necessary modules" "Imported the necessary modules/libraries"

# 3. "date different of the START date" -> "date difference of the START date"
Replace-Text "date different of the START date" "date difference of the START date"

# 4. "73 for patient " -> "73 for patients "
Replace-Text "73 for patient " "73 for patients "

# 5. "used One Hot Encoder to transform" -> "Used One Hot Encoder to transform"
Replace-Text "used One Hot Encoder to transform" "Used One Hot Encoder to transform"

# 6. ", with rows of 75% and 25% respectively, rows are selected randomly using scikit-learn " ->
#    ", with rows of 75% and 25% respectively. Rows are selected randomly using scikit-learn "
Replace-Text ", with rows of 75% and 25% respectively, rows are selected randomly using scikit-learn " ", with rows of 75% and 25% respectively. Rows are selected randomly using scikit-learn "

# 7. "coefficients show that the AGE is important feature, next most important feature is GENDER" ->
#    "coefficients show that the AGE is an important feature, next most important feature is GENDER"
Replace-Text "coefficients show that the AGE is important feature, next most important feature is GENDER" "coefficients show that the AGE is an important feature, next most important feature is GENDER"

# 8. " the rest of modeling, RACE was kept" -> " in the rest of modeling process, RACE was kept"
Replace-Text " the rest of modeling, RACE was kept" " in the rest of modeling process, RACE was kept"

# 9. (run split only, no text change - skipped)

# 10. "hidden network layers, the number of nodes in a layer" ->
#     "hidden network layers, and/or the number of nodes in a layer"
Replace-Text "hidden network layers, the number of nodes in a layer" "hidden network layers, and/or the number of nodes in a layer"

# 11. " the model creation (create.py), does the similar things as codes in " ->
#     " The model creation (create.py), which does the similar process as the codes do in "
Replace-Text " the model creation (create.py), does the similar things as codes in " " The model creation (create.py), which does the similar process as the codes do in "

# 12. "GUI, an application (predict.py), provides preliminary graphic user interface, the user " ->
#     "GUI, an application (predict.py), provides preliminary graphic user interface, so that the user "
Replace-Text "GUI, an application (predict.py), provides preliminary graphic user interface, the user " "GUI, an application (predict.py), provides preliminary graphic user interface, so that the user "

# 13. " and enter the age, the click the Predict button to see the results from 4 classic models" ->
#     " and enter the age, then click the Predict button to see the results from 4 classic models (Survive or Die)"
Replace-Text " and enter the age, the click the Predict button to see the results from 4 classic models" " and enter the age, then click the Predict button to see the results from 4 classic models (Survive or Die)"

# 14. "It is in private mode, but access is granted to individual by invite" ->
#     "It is in private mode, but access can be granted to individual(s) by invite"
Replace-Text "It is in private mode, but access is granted to individual by invite" "It is in private mode, but access can be granted to individual(s) by invite"
